# Fix dataset problems in "commenti_positivi_negativi_per_notizia_cronaca_nera"
# Corrects negativo/positivo counts and giornale/social labels for several
# rows (CRONACA NERA topic) and adds a missing data row (La Repubblica /
# YouTube / Strage di Cutro) that had been dropped from the original export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: row, topic(B), titolo(C), giornale(D), social(E), negativo(F), positivo(G)
$rows = @(
  @(3,  'CRONACA NERA', 'Incidente Youtubers', 'FanPage', 'Facebook', 89, 11),
  @(4,  'CRONACA NERA', 'Incidente Youtubers', 'FanPage', 'Instagram', 98, 2),
  @(5,  'CRONACA NERA', 'Incidente Youtubers', 'FanPage', 'YouTube', 90, 10),
  @(6,  'CRONACA NERA', 'Incidente Youtubers', 'Il Corriere Della Sera', 'Facebook', 94, 6),
  @(7,  'CRONACA NERA', 'Incidente Youtubers', 'Il Corriere Della Sera', 'Instagram', 94, 6),
  @(8,  'CRONACA NERA', 'Incidente Youtubers', 'Il Corriere Della Sera', 'YouTube', 94, 6),
  @(9,  'CRONACA NERA', 'Incidente Youtubers', 'La Repubblica', 'Facebook', 86, 14),
  @(10, 'CRONACA NERA', 'Incidente Youtubers', 'La Repubblica', 'Instagram', 91, 9),
  @(11, 'CRONACA NERA', 'Incidente Youtubers', 'La Repubblica', 'YouTube', 93, 7),
  @(12, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'FanPage', 'Facebook', 88, 12),
  @(13, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'FanPage', 'Instagram', 76, 24),
  @(14, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'FanPage', 'YouTube', 77, 23),
  @(15, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'Il Corriere Della Sera', 'Facebook', 89, 11),
  @(16, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'Il Corriere Della Sera', 'Instagram', 83, 17),
  @(17, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'Il Corriere Della Sera', 'YouTube', 84, 16),
  @(18, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'La Repubblica', 'Facebook', 85, 15),
  @(19, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'La Repubblica', 'Instagram', 87, 13),
  @(20, 'CRONACA NERA', "L'implosione del sottomarino Titan", 'La Repubblica', 'YouTube', 79, 21),
  @(21, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'FanPage', 'Facebook', 79, 21),
  @(22, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'FanPage', 'Instagram', 88, 12),
  @(23, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'FanPage', 'YouTube', 79, 21),
  @(24, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'Il Corriere Della Sera', 'Facebook', 91, 9),
  @(25, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'Il Corriere Della Sera', 'Instagram', 78, 22),
  @(26, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'Il Corriere Della Sera', 'YouTube', 88, 12),
  @(27, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'La Repubblica', 'Facebook', 87, 13),
  @(28, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'La Repubblica', 'Instagram', 87, 13),
  @(29, 'CRONACA NERA', "L'omicidio di Giulia Cecchettin", 'La Repubblica', 'YouTube', 82, 18),
  @(30, 'CRONACA NERA', 'Strage di Cutro', 'FanPage', 'Facebook', 93, 7),
  @(31, 'CRONACA NERA', 'Strage di Cutro', 'FanPage', 'Instagram', 91, 9),
  @(32, 'CRONACA NERA', 'Strage di Cutro', 'FanPage', 'YouTube', 87, 13),
  @(33, 'CRONACA NERA', 'Strage di Cutro', 'Il Corriere Della Sera', 'Facebook', 90, 10),
  @(34, 'CRONACA NERA', 'Strage di Cutro', 'Il Corriere Della Sera', 'Instagram', 92, 8),
  @(35, 'CRONACA NERA', 'Strage di Cutro', 'Il Corriere Della Sera', 'YouTube', 88, 12),
  @(36, 'CRONACA NERA', 'Strage di Cutro', 'La Repubblica', 'Facebook', 93, 7),
  @(37, 'CRONACA NERA', 'Strage di Cutro', 'La Repubblica', 'Instagram', 79, 21),
  @(38, 'CRONACA NERA', 'Strage di Cutro', 'La Repubblica', 'YouTube', 82, 18)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
}
